# Applies the "Tested different solutions for setting the worker values." edit
# to the schedule workbook: updates the per-day shift assignments (columns F:L)
# for the worker rows, appends a "test" value, removes a stray value, and
# updates the sheet view selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update / add shift-assignment cells (columns F:L) ---------------------
$ws.Range("F2").Value  = "07-15;"
$ws.Range("G2").Value  = "07-15;"

$ws.Range("F3").Value  = "07-15;"
$ws.Range("G3").Value  = "15-23;"

$ws.Range("F4").Value  = "07-15;"
$ws.Range("G4").Value  = "23-07;"

$ws.Range("F5").Value  = "07-15;15-23;"
$ws.Range("G5").Value  = "07-15;15-23;"
$ws.Range("H5").Value  = "07-15;15-23;"

$ws.Range("F6").Value  = "23-07;"
$ws.Range("G6").Value  = "23-07;"
$ws.Range("H6").Value  = "23-07;"
$ws.Range("I6").Value  = "23-07;"
$ws.Range("J6").Value  = "23-07;"
$ws.Range("K6").Value  = "23-07;"

$ws.Range("F7").Value  = "15-23;23-07;"
$ws.Range("G7").Value  = "15-23;23-07;"
$ws.Range("H7").Value  = "15-23;23-07;"

$ws.Range("F8").Value  = "15-23;"
$ws.Range("G8").Value  = "15-23;"
$ws.Range("H8").Value  = "15-23;"
$ws.Range("I8").Value  = "15-23;"
$ws.Range("J8").Value  = "15-23;"
$ws.Range("K8").Value  = "15-23;"
$ws.Range("L8").Value  = "07-15;15-23;23-07;"

$ws.Range("F9").Value  = "15-23;23-07;"
$ws.Range("H9").Value  = "07-15;15-23;23-07;"
$ws.Range("I9").Value  = "07-15;"
$ws.Range("J9").Value  = "07-15;15-23;23-07;"
$ws.Range("K9").Value  = "07-15;15-23;23-07;"
$ws.Range("L9").Value  = "07-15;15-23;23-07;"

$ws.Range("K10").Value = "07-15;15-23;23-07;"
$ws.Range("L10").Value = "07-15;15-23;23-07;"

$ws.Range("G11").Value = "07-15;15-23;"
$ws.Range("I11").Value = "15-23;"
$ws.Range("J11").Value = "15-23;"
$ws.Range("L11").Value = "15-23;"

$ws.Range("I12").Value = "15-23;"
$ws.Range("J12").Value = "07-15;"
$ws.Range("K12").Value = "07-15;23-07;"

$ws.Range("H13").Value = "15-23;23-07;"
$ws.Range("I13").Value = "07-15;"
$ws.Range("J13").Value = "23-07;"
$ws.Range("K13").Value = "07-15;15-23;"

$ws.Range("F14").Value = "15-23;"
$ws.Range("G14").Value = "07-15;"
$ws.Range("H14").Value = "23-07;"
$ws.Range("I14").Value = "15-23;"
$ws.Range("K14").Value = "23-07;"

$ws.Range("F15").Value = "15-23;"
$ws.Range("G15").Value = "07-15;15-23;"
$ws.Range("I15").Value = "07-15;"
$ws.Range("J15").Value = "07-15;15-23;"

$ws.Range("F16").Value = "07-15;"
$ws.Range("G16").Value = "15-23;"
$ws.Range("I16").Value = "15-23;"
$ws.Range("J16").Value = "23-07;"
$ws.Range("K16").Value = "23-07;"
$ws.Range("L16").Value = "07-15;"

$ws.Range("H17").Value = "07-15;15-23;23-07;"
$ws.Range("J17").Value = "15-23;"
$ws.Range("K17").Value = "07-15;"
$ws.Range("L17").Value = "07-15;"

$ws.Range("F18").Value = "07-15;"
$ws.Range("I18").Value = "07-15;"
$ws.Range("J18").Value = "23-07;"
$ws.Range("K18").Value = "15-23;23-07;"

$ws.Range("G19").Value = "15-23;"
$ws.Range("H19").Value = "15-23;"

$ws.Range("F20").Value = "07-15;15-23;23-07;"
$ws.Range("K20").Value = "07-15;15-23;23-07;"

$ws.Range("G21").Value = "07-15;15-23;23-07;"

$ws.Range("K22").Value = "07-15;15-23;23-07;"

$ws.Range("F23").Value = "15-23;"
$ws.Range("G23").Value = "15-23;"
$ws.Range("I23").Value = "07-15;23-07;"
$ws.Range("K23").Value = "test"
$ws.Range("L23").Value = "07-15;23-07;"

$ws.Range("F24").Value = "15-23;23-07;"
$ws.Range("G24").Value = "07-15;15-23;"
$ws.Range("I24").Value = "15-23;"
$ws.Range("J24").Value = "15-23;"

$ws.Range("F25").Value = "15-23;"
$ws.Range("G25").Value = "15-23;"
$ws.Range("H25").Value = "15-23;"
$ws.Range("I25").Value = "15-23;"
$ws.Range("J25").Value = "15-23;"
$ws.Range("K25").Value = "15-23;"

# L25 no longer holds a value ("Test" was removed).
$ws.Range("L25").Clear()

# --- Update the sheet view: scroll position + active selection -------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("L15").Select()
